# Auto-generated Excel COM-interop edit script
# Applies the weekly NYPD CompStat crime-data refresh to the CompStat_1 sheet:
#   - bumps the report Volume/Number and the covered date range
#   - refreshes all weekly/28-day/YTD crime counts and percentage changes
#     for precincts/rows 16-30

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# Header text updates (Volume number, report date range).
# Use Characters() so only the relevant substring of the rich-text
# cell is touched.
# ------------------------------------------------------------------
$volCell = $ws.Range("A8")
$volCell.Characters(21, 2).Text = "51"

$dateCell = $ws.Range("C9")
$dateCell.Characters(27, 10).Text = "12/18/2023"
$dateCell.Characters(48, 10).Text = "12/24/2023"

# ------------------------------------------------------------------
# Donor cells used purely to copy a specific, already-existing cell
# style (number formats) onto cells whose value "type" changes below
# (shared-string placeholder <-> real number).
# ------------------------------------------------------------------
$textStyleDonor = $ws.Range("C14")   # general/text style for "---" placeholders

# ------------------------------------------------------------------
# Crime-statistics table updates (rows 16-30)
# ------------------------------------------------------------------
$ws.Range("C16").Value = 2
$ws.Range("D16").Value = 1
$ws.Range("D16").NumberFormat = '#,##0'
$ws.Range("E16").Value = 100
$ws.Range("E16").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("F16").Value = 18
$ws.Range("G16").Value = 4
$ws.Range("H16").Value = 350
$ws.Range("I16").Value = 135
$ws.Range("J16").Value = 149
$ws.Range("K16").Value = -9.395973154362
$ws.Range("L16").Value = -4.929577464788
$ws.Range("M16").Value = -28.571428571428
$ws.Range("N16").Value = -84.536082474226
$ws.Range("D17").Value = 6
$ws.Range("E17").Value = -66.666666666666
$ws.Range("F17").Value = 9
$ws.Range("H17").Value = -25
$ws.Range("I17").Value = 176
$ws.Range("J17").Value = 186
$ws.Range("K17").Value = -5.376344086021
$ws.Range("L17").Value = -3.825136612021
$ws.Range("M17").Value = 76
$ws.Range("N17").Value = -61.487964989059
$ws.Range("C18").Value = 3
$ws.Range("D18").Value = 3
$ws.Range("E18").Value = 0
$ws.Range("F18").Value = 7
$ws.Range("G18").Value = 14
$ws.Range("H18").Value = -50
$ws.Range("I18").Value = 168
$ws.Range("J18").Value = 180
$ws.Range("K18").Value = -6.666666666666
$ws.Range("L18").Value = 27.272727272727
$ws.Range("M18").Value = 34.4
$ws.Range("N18").Value = -84.060721062618
$ws.Range("C19").Value = 7
$ws.Range("D19").Value = 8
$ws.Range("E19").Value = -12.5
$ws.Range("F19").Value = 30
$ws.Range("G19").Value = 38
$ws.Range("H19").Value = -21.052631578947
$ws.Range("I19").Value = 501
$ws.Range("J19").Value = 520
$ws.Range("K19").Value = -3.653846153846
$ws.Range("L19").Value = 7.051282051282
$ws.Range("M19").Value = -5.471698113207
$ws.Range("N19").Value = -54.660633484162
$ws.Range("D20").Value = 2
$ws.Range("G20").Value = 9
$ws.Range("H20").Value = -44.444444444444
$ws.Range("J20").Value = 88
$ws.Range("K20").Value = 7.954545454545
$ws.Range("L20").Value = 7.954545454545
$ws.Range("M20").Value = 115.909090909091
$ws.Range("N20").Value = -90.306122448979
$ws.Range("C21").Value = 14
$ws.Range("D21").Value = 20
$ws.Range("E21").Value = -30
$ws.Range("F21").Value = 71
$ws.Range("G21").Value = 77
$ws.Range("H21").Value = -7.792207792207
$ws.Range("I21").Value = 1089
$ws.Range("J21").Value = 1133
$ws.Range("K21").Value = -3.883495145631
$ws.Range("L21").Value = 5.933852140077
$ws.Range("M21").Value = 8.791208791208
$ws.Range("N21").Value = -75.944333996023
$ws.Range("C22").Value = 1
$ws.Range("C22").NumberFormat = '#,##0'
$ws.Range("D22").Value = 1
$ws.Range("D22").NumberFormat = '#,##0'
$ws.Range("E22").Value = 0
$ws.Range("E22").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("F22").Value = 5
$ws.Range("G22").Value = 4
$ws.Range("H22").Value = 25
$ws.Range("I22").Value = 26
$ws.Range("J22").Value = 33
$ws.Range("K22").Value = -21.212121212121
$ws.Range("L22").Value = -7.142857142857
$ws.Range("M22").Value = -13.333333333333
$ws.Range("C23").Value = "'0"
$textStyleDonor.Copy() | Out-Null
$ws.Range("C23").PasteSpecial(-4122)
$ws.Range("E23").Value = -100
$ws.Range("F23").Value = 3
$ws.Range("G23").Value = 4
$ws.Range("H23").Value = -25
$ws.Range("J23").Value = 121
$ws.Range("K23").Value = -11.570247933884
$ws.Range("L23").Value = 0
$ws.Range("M23").Value = 42.666666666666
$ws.Range("C24").Value = 35
$ws.Range("D24").Value = 32
$ws.Range("E24").Value = 9.375
$ws.Range("F24").Value = 104
$ws.Range("G24").Value = 118
$ws.Range("H24").Value = -11.864406779661
$ws.Range("I24").Value = 1537
$ws.Range("J24").Value = 1821
$ws.Range("K24").Value = -15.595826468973
$ws.Range("L24").Value = 18.962848297213
$ws.Range("M24").Value = 45.963912630579
$ws.Range("C25").Value = 3
$ws.Range("D25").Value = 3
$ws.Range("E25").Value = 0
$ws.Range("F25").Value = 21
$ws.Range("H25").Value = 10.526315789473
$ws.Range("I25").Value = 300
$ws.Range("J25").Value = 299
$ws.Range("K25").Value = 0.334448160535
$ws.Range("L25").Value = 10.294117647058
$ws.Range("M25").Value = -6.542056074766
$ws.Range("C27").Value = "'0"
$textStyleDonor.Copy() | Out-Null
$ws.Range("C27").PasteSpecial(-4122)
$ws.Range("D27").Value = 3
$ws.Range("D27").NumberFormat = '#,##0'
$ws.Range("E27").Value = -100
$ws.Range("E27").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("F27").Value = 2
$ws.Range("G27").Value = 3
$ws.Range("H27").Value = -33.333333333333
$ws.Range("I27").Value = 42
$ws.Range("J27").Value = 54
$ws.Range("K27").Value = -22.222222222222
$ws.Range("L27").Value = -19.230769230769
$ws.Range("N28").Value = -87.272727272727
$ws.Range("N29").Value = -88.461538461538
$ws.Range("F30").Value = 1

$excel.CutCopyMode = $false

